$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the new B:J values for one row (rows 2-25).
$data = @(
    @(4.377150980694014, 3.379914246651969, 9.449768375945162, 18.47399760305157, 7.231757402547101, 2.960181145736489, 2.507522494578346, 3.414291668632348, 1.476604685401425),
    @(2.751815514474119, 4.284371319845489, 10.1342376280923, 18.52376919492934, 6.767798014299466, 2.960251473765107, 2.614411001142907, 3.454475047340116, 1.618687900100354),
    @(2.757074984518666, 4.259200753344929, 9.718818166843116, 18.48597518109716, 6.397200910192439, 2.960311570509332, 2.70568851398221, 3.409662268348357, 1.515560131842274),
    @(2.763619456296708, 4.315873001313871, 9.581627505098865, 18.37166118868871, 6.20059484244469, 2.960342804261893, 2.753176120253797, 3.292756276489932, 1.617321895261638),
    @(2.790119313091753, 4.419894277330862, 9.615804525064506, 18.28610256189373, 6.114846664261417, 2.960355899475239, 2.77307833132996, 3.206609321741567, 1.643792517330155),
    @(3.457248895295045, 5.193241574679397, 10.50432396482842, 18.2858021046148, 6.114857602873973, 2.960350680735095, 2.765145861482968, 3.206609321741611, 1.698448044258332),
    @(3.733501088231413, 5.397867408501163, 10.75497804858412, 18.26904221803913, 6.136077023567042, 2.960345865583092, 2.757787240833748, 3.190568529478301, 1.812755433389108),
    @(2.98048710186955, 4.972173630366306, 10.59690162634624, 18.23808155454232, 6.442051925991617, 2.960299965399589, 2.688024012162442, 3.165909919073763, 1.886585556293671),
    @(4.494599533892234, 3.320709194791647, 9.329830122772751, 18.30524271684992, 7.197565067547022, 2.960187414781246, 2.516972042956099, 3.246804273437029, 1.489390065793366),
    @(9.329644214908447, 2.428127632827539, 7.808595192214611, 18.34788288210802, 7.917319569547963, 2.960074358164766, 2.345059120591833, 3.303339636634402, 1.020937903443516),
    @(11.14834013443572, 2.557384939984241, 7.320617970422393, 18.3155584387526, 8.258310533330446, 2.960017830615906, 2.273604417646808, 3.278512933571507, 0.8988582028775227),
    @(11.65151335541934, 2.599030371490857, 6.93684330454513, 18.28188141503194, 8.522244728219166, 2.959973387653754, 2.372036402585024, 3.25084277283168, 0.8027070478866629),
    @(11.05542827332309, 2.497707232187481, 6.835572182382715, 18.24054239096546, 8.568401783898235, 2.959965675474017, 2.595277714579796, 3.21104898746755, 0.7861683681731),
    @(10.42371651852407, 2.442104697467491, 6.866949929700393, 18.18205877504704, 8.501862865039907, 2.959977087982375, 2.722068456971865, 3.151996471117923, 0.8113920867200061),
    @(10.1747940245792, 2.428030667570541, 6.922923344991936, 18.24060236535438, 8.502499199550655, 2.959976763736356, 2.765078240612815, 3.20971201073794, 0.8104340253762247),
    @(9.946360435311478, 2.420897897676972, 6.991290687205554, 18.31654775833172, 8.506557338280803, 2.959975800024422, 2.796948647756782, 3.284725098375376, 0.8079992017358351),
    @(8.840329375667626, 2.45264855696137, 7.275736778961825, 18.31381049561193, 8.2906262830847, 2.960011966157307, 2.815457430931946, 3.277477626738684, 0.8869254641439041),
    @(8.350794688711122, 2.491750244930937, 7.490444667626228, 18.35066629042786, 8.156694602634566, 2.960034132227903, 2.74589474052704, 3.311069770235996, 0.9349190441063394),
    @(9.540798834865871, 2.424662992557593, 7.482688034052917, 18.34398180812121, 8.157681858686605, 2.96003443778853, 2.45000274239669, 3.304476368388815, 0.9348594902234381),
    @(12.49912417773591, 2.899309808985555, 7.047593734445005, 18.37502472245938, 8.506950527730137, 2.959975549455284, 2.203215430056046, 3.342478804598323, 0.8071408605508809),
    @(14.35414606553618, 3.508940837538074, 6.676126401232599, 18.45452127466651, 8.845858291679059, 2.959916829024916, 2.106081306190618, 3.428282708650405, 0.6798576488820666),
    @(15.65611241077847, 3.9615742293201, 6.419946305125531, 18.54252701902952, 9.105157622292944, 2.959871059564684, 2.036416869975091, 3.520942166688997, 0.5803980712816041),
    @(16.07935722637156, 4.107669621934711, 6.355059738335359, 18.60632125471107, 9.201526196288873, 2.959853833222639, 2.010135366830154, 3.586192899609446, 0.5426936739556507),
    @(13.09332025580006, 3.057856359654183, 7.078035576335987, 18.66584580469569, 8.697019707373684, 2.959942196474851, 2.144414751122814, 3.634292947091492, 0.7331531971415753),
)

$startRow = 2
$startCol = 2  # column B

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowValues[$j]
    }
}
